# Assignment 1.docx - "Started designing the BTPAdministratorClient class"
#
# The author clicked inside the requirements table, placing the text
# cursor between "Apply an" and "interest rate on a savings account"
# (right before the separating space). Word records the position of the
# most recent edit with its reserved "_GoBack" bookmark, which forces the
# run containing that text to split around the (now zero-length) bookmark.
# Because a document may only have a single "_GoBack" bookmark, adding the
# new one automatically retires the "_GoBack" bookmark that used to sit
# later in the document (just before "the amount to "), and the other,
# user-named bookmarks that sorted after it are renumbered accordingly by
# Word when it re-serialises the bookmark ids.

$d = $word.ActiveDocument

# The row "11 | Apply an interest rate on a savings account | -" lives in
# the 3rd table in the document (row 8, column 2).
$table = $d.Tables.Item(3)
$cell = $table.Cell(8, 2)
$cellStart = $cell.Range.Start

# "Apply an interest rate on a savings account"
#  0123456789...
# The space that separates "Apply an" from "interest..." is at offset 8
# (0-based) from the start of the cell.
$splitOffset = 8

# 1) Drop the "_GoBack" bookmark right after "Apply an" - this both marks
#    the edit point and splits "Apply an" off from the remaining text.
$splitPos = $cellStart + $splitOffset
$bookmarkRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# 2) Separate "interest rate on a savings account" into its own run by
#    briefly annotating it with a comment, then discarding the comment -
#    this leaves the run boundary behind without leaving any formatting
#    residue.
$interestStart = $cellStart + $splitOffset + 1
$interestEnd = $cell.Range.End - 1
$interestRange = $d.Range($interestStart, $interestEnd)
$tempComment = $d.Comments.Add($interestRange, "split")
$tempComment.Delete()
